$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block updates ---------------------------------------------
# Row 2: MatrNr becomes a real number instead of the placeholder "?"
$ws.Range("B2").Value = 3139954

# Row 4: Team name becomes "Rayforge" instead of the generic project name
$ws.Range("B4").Value = "Rayforge"

# --- New time-tracking rows (8-18) -------------------------------------
# Data: row, date-serial, hours, Buchungsposten, Beschreibung
$entries = @(
    @(8,  45579, 1.5, "Online-Meeting",     "Projekt Kickoff Meeting"),
    @(9,  45580, 1,   "Recherche",          "Bezüglich .obj / .mtl Dateien"),
    @(10, 45583, 2,   "Recherche",          "Bezüglich .obj / .mtl Dateien"),
    @(11, 45586, 1.5, "Online-Meeting",     "Weekly-Summup-01 Meeting"),
    @(12, 45586, 1.5, "Coding",             "Erste Ansätze im Java Object Parser"),
    @(13, 45587, 1,   "Emergency-Meeting",  "Emergency-Meeting abgehalten. (Mail von Management missverstanden)"),
    @(14, 45589, 1.5, "Coding",             "Umschreiben des Object Parsers"),
    @(15, 45591, 3,   "Recherche",          "Bezüglich .obj / .mtl Dateien und erstellen einer Powerpoint für Teammitglieder"),
    @(16, 45592, 2,   "Coding",             "Fertigstellen des Object Parser Prototyps"),
    @(17, 45593, 2,   "Online-Meeting",     "Präsentation meiner .obj / .mtl Informationen und meines Prototyps"),
    @(18, 45593, 1.5, "Recherche",          "Auseinandersetzen mit neuem OpenGL Code")
)

foreach ($entry in $entries) {
    $row  = $entry[0]
    $date = $entry[1]
    $hrs  = $entry[2]
    $cat  = $entry[3]
    $desc = $entry[4]

    # Carry the date/number formatting down from the template row (7)
    # so the new rows reuse the existing cell styles instead of creating new ones.
    $ws.Range("A7").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("B7").Copy() | Out-Null
    $ws.Range("B$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$row").Value = $date
    $ws.Range("B$row").Value = $hrs
    $ws.Range("C$row").Value = $cat
    $ws.Range("D$row").Value = $desc
}

$excel.CutCopyMode = 0

# --- View state (zoom level & selected cell) ---------------------------
$excel.ActiveWindow.Zoom = 144
$ws.Range("D21").Select() | Out-Null
